# Update gh-pages output - refresh "想去人数" (F column) counts across sheets
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 181
$ws1.Range("F3").Value = 397
$ws1.Range("F4").Value = 1114
$ws1.Range("F5").Value = 35
$ws1.Range("F7").Value = 13
$ws1.Range("F8").Value = 1052
$ws1.Range("F10").Value = 312
$ws1.Range("F11").Value = 408
$ws1.Range("F13").Value = 298
$ws1.Range("F15").Value = 20
$ws1.Range("F17").Value = 382
$ws1.Range("F18").Value = 430
$ws1.Range("F19").Value = 5469
$ws1.Range("F20").Value = 86
$ws1.Range("F21").Value = 1532
$ws1.Range("F22").Value = 349
$ws1.Range("F23").Value = 4608
$ws1.Range("F24").Value = 4608
$ws1.Range("F27").Value = 1465
$ws1.Range("F28").Value = 9
$ws1.Range("F30").Value = 635
$ws1.Range("F31").Value = 16
$ws1.Range("F33").Value = 3782

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 101

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 9368
$ws3.Range("F3").Value = 578
$ws3.Range("F4").Value = 2111

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9368
$ws4.Range("F3").Value = 578
$ws4.Range("F4").Value = 2111
$ws4.Range("F5").Value = 181
$ws4.Range("F6").Value = 397
$ws4.Range("F7").Value = 1114
$ws4.Range("F8").Value = 35
$ws4.Range("F10").Value = 13
$ws4.Range("F11").Value = 1052
$ws4.Range("F13").Value = 312
$ws4.Range("F14").Value = 408
$ws4.Range("F16").Value = 298
$ws4.Range("F18").Value = 20
$ws4.Range("F23").Value = 382
$ws4.Range("F24").Value = 430
$ws4.Range("F25").Value = 5469
$ws4.Range("F26").Value = 86
$ws4.Range("F27").Value = 1532
$ws4.Range("F30").Value = 349
$ws4.Range("F32").Value = 4608
$ws4.Range("F33").Value = 4608
$ws4.Range("F36").Value = 1465
$ws4.Range("F37").Value = 9
$ws4.Range("F39").Value = 635
$ws4.Range("F40").Value = 16
$ws4.Range("F47").Value = 3783
